$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly column AB: date header "11_05_2021", values for rows 2-11, and
# the SUM formula in row 12 (extending the existing shared-formula range).
$ws.Range("AB1").Value = "11_05_2021"

$ws.Range("AB2").Value = 206
$ws.Range("AB3").Value = 209
$ws.Range("AB4").Value = 675
$ws.Range("AB5").Value = 1034
$ws.Range("AB6").Value = 1483
$ws.Range("AB7").Value = 2321
$ws.Range("AB8").Value = 2331
$ws.Range("AB9").Value = 3207
$ws.Range("AB10").Value = 2559
$ws.Range("AB11").Value = 704

$ws.Range("AB12").Formula = "=SUM(AB2:AB11)"

# Match the saved view state from the diff (scrolled / selected range moved
# one column right to keep the newest week in view).
$ws.Range("AB14").Select()
$excel.ActiveWindow.ScrollColumn = 19
